# Append two new rows (9 and 10) to Sheet1, duplicating rows 7 and 8
# (the latest "a1"/"b2" samples), as new data points starting from the
# 980 initial-value recommendation feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceRows = @(7, 8)
$destRows = @(9, 10)

for ($i = 0; $i -lt $sourceRows.Length; $i++) {
    $srcRow = $sourceRows[$i]
    $dstRow = $destRows[$i]

    # Copy number/cell formatting (borders, font, alignment) from the
    # source row onto the new row first.
    $srcRowRange = $ws.Range("A$srcRow`:K$srcRow")
    $dstRowRange = $ws.Range("A$dstRow`:K$dstRow")
    $srcRowRange.Copy()
    $dstRowRange.PasteSpecial(-4122)  # xlPasteFormats

    # Then copy the actual cell values across (Value2 avoids the
    # reflection-artifact bug seen when reading/writing .Value directly).
    for ($col = 1; $col -le 11; $col++) {
        $srcCell = $ws.Cells.Item($srcRow, $col)
        $dstCell = $ws.Cells.Item($dstRow, $col)
        $dstCell.Value2 = $srcCell.Value2
    }
}

$excel.CutCopyMode = $false
